$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before DM (column 117), shifting DM:DN onward (nom,
# url_produit, ...) one column to the right, i.e. to DN:DO. This matches the
# target dimension change A1:DN206 -> A1:DO206. The newly inserted column
# inherits formatting from its left neighbour (DL), so the bold/centered
# header style (s="1") carries over automatically.
$ws.Range("DM1").EntireColumn.Insert()

# New timestamp header for the freshly inserted snapshot column.
$ws.Range("DM1").Value = "2026-02-02 01:02:19"

# For every data row that still has a recorded price history (rows 2-80),
# duplicate the last recorded price (column DL = 116, untouched by the
# insert) into the newly inserted column DM (= 117) so the new snapshot
# carries the last known price forward.
for ($r = 2; $r -le 80; $r++) {
    $ws.Cells.Item($r, 117).Value = $ws.Cells.Item($r, 116).Value()
}
